# Generate Report for Handback
#
# Mirrors the localization-status report refresh: the zh-cn / de-de handback
# round-trip completed, so:
#   - status text moves from "Ready for handoff" to "Handed back: in sync
#     with en-US" everywhere it is shown (Overview + per-locale sheets),
#   - each locale sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get filled in with the handback
#     artifact names + timestamp, with the target-file cell turned into a
#     hyperlink (matching the existing A-column hyperlink look),
#   - a few columns widen to fit the longer values that now live in them.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) == the workbook's existing "FF6495ED" hyperlink font color

function Set-AsHyperlink($ws, $cellAddr, $displayText, $url) {
    $range = $ws.Range($cellAddr)
    $ws.Hyperlinks.Add($range, $url, "", "", $displayText) | Out-Null
    # Hyperlinks.Add applies Excel's built-in theme-colored "Hyperlink" cell
    # style; restore the look already used elsewhere in this workbook for
    # hyperlinked cells (underlined, custom blue font color).
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: the per-locale status cells flip to the new wording.
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $statusNew
$ovw.Range("F2").Value = $statusNew
$ovw.Range("E3").Value = $statusNew
$ovw.Range("F3").Value = $statusNew

$ovw.Columns.Item(5).ColumnWidth = 29.9777047293527
$ovw.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

Set-AsHyperlink $zh "I2" "1e66d068-0e87-4431-917f-2ea2c4a17a11.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c37dee6dcbac1cc172b6f366192cb264c97318d/e2e/1e66d068-0e87-4431-917f-2ea2c4a17a11.md"
$zh.Range("J2").Value = "1e66d068-0e87-4431-917f-2ea2c4a17a11.6ae6942611cb798255de15786bef033def6451f0.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-04 16:33:41"

Set-AsHyperlink $zh "I3" "95cb353e-a87a-4d65-92e0-9fdfefb4090e.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c37dee6dcbac1cc172b6f366192cb264c97318d/e2e/95cb353e-a87a-4d65-92e0-9fdfefb4090e.md"
$zh.Range("J3").Value = "95cb353e-a87a-4d65-92e0-9fdfefb4090e.03a0296777bed3b048df9fdb0a78591c7b9b6f47.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-04 16:33:41"

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

Set-AsHyperlink $de "I2" "1e66d068-0e87-4431-917f-2ea2c4a17a11.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c37dee6dcbac1cc172b6f366192cb264c97318d/e2e/1e66d068-0e87-4431-917f-2ea2c4a17a11.md"
$de.Range("J2").Value = "1e66d068-0e87-4431-917f-2ea2c4a17a11.6ae6942611cb798255de15786bef033def6451f0.de-de.xlf"
$de.Range("K2").Value = "2016-09-04 16:33:48"

Set-AsHyperlink $de "I3" "95cb353e-a87a-4d65-92e0-9fdfefb4090e.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c37dee6dcbac1cc172b6f366192cb264c97318d/e2e/95cb353e-a87a-4d65-92e0-9fdfefb4090e.md"
$de.Range("J3").Value = "95cb353e-a87a-4d65-92e0-9fdfefb4090e.03a0296777bed3b048df9fdb0a78591c7b9b6f47.de-de.xlf"
$de.Range("K3").Value = "2016-09-04 16:33:48"

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40

Write-Host "Handback report generated."
